$wb = $excel.ActiveWorkbook

# Remember which sheet is active/selected so we can restore it at the end -
# adding a worksheet below implicitly activates the new sheet.
$originallyActiveName = $wb.ActiveSheet.Name

# ---------------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet - insert a new row for 2022-Q3
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.29

# A2 should carry the same style as the other "序号" cells below it (A3:A6)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: insert a brand-new worksheet "2022-Q3" right after "总计"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Copy the header-row + first-column formatting from the "2022-Q2" sheet
# (now shifted to position 3) so the new sheet matches the existing layout.
$refSheet = $wb.Worksheets.Item(3)

$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 3: populate the header row
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# Step 4: populate the data rows (2 - 5)
#   Column A & H are real numbers; columns B-G must stay TEXT (fund codes
#   keep leading zeros, and the decimal-looking figures are stored as text
#   in the original workbook) so force Text format before writing, then
#   clear the formatting back to the default (unstyled) cell style.
# ---------------------------------------------------------------------------
$textRange = $newSheet.Range("B2:G5")
$textRange.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "005585"
$newSheet.Range("C2").Value = "银河文体娱乐主题灵活配置混合A"
$newSheet.Range("D2").Value = "3.01"
$newSheet.Range("E2").Value = "90.28"
$newSheet.Range("F2").Value = "4.64"
$newSheet.Range("G2").Value = "0.1397"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "001628"
$newSheet.Range("C3").Value = "招商体育文化休闲股票A"
$newSheet.Range("D3").Value = "2.23"
$newSheet.Range("E3").Value = "92.42"
$newSheet.Range("F3").Value = "5.20"
$newSheet.Range("G3").Value = "0.1160"
$newSheet.Range("H3").Value = 2

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "015667"
$newSheet.Range("C4").Value = "银河文体娱乐主题灵活配置混合C"
$newSheet.Range("D4").Value = "0.41"
$newSheet.Range("E4").Value = "90.28"
$newSheet.Range("F4").Value = "4.64"
$newSheet.Range("G4").Value = "0.0190"
$newSheet.Range("H4").Value = 7

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "015395"
$newSheet.Range("C5").Value = "招商体育文化休闲股票C"
$newSheet.Range("D5").Value = "0.25"
$newSheet.Range("E5").Value = "92.42"
$newSheet.Range("F5").Value = "5.20"
$newSheet.Range("G5").Value = "0.0130"
$newSheet.Range("H5").Value = 2

$textRange.ClearFormats()

# Re-apply the copied formatting (ClearFormats above wiped it) for the
# header row and column A so the layout matches the other quarter sheets.
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 5: restore the originally active/selected sheet (adding $newSheet
# above made it the active sheet as a side effect).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($originallyActiveName).Activate()

Write-Output "done"
